$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 78.260000000000005
$ws.Range("D2").Value = 69

$ws.Range("B3").Value = 70.83
$ws.Range("D3").Value = 48

$ws.Range("B4").Value = 72.22
$ws.Range("D4").Value = 72

$ws.Range("C5").Value = 86.09
$ws.Range("E5").Value = 115

$ws.Range("C6").Value = 86.67
$ws.Range("E6").Value = 135

$ws.Range("C7").Value = 83.23
$ws.Range("E7").Value = 155

$ws.Range("E8").Select()
